# Guardar cambios: solicitud de usuario
# Restructure the "Clientes" sheet: rename/replace several header columns
# (F..V), append new trailing header columns (W..AC) that carry forward
# fecha_ingreso/fecha_dispersion/fecha_proximo/estatus/observaciones/score/
# sucursal, and relocate the matching row-2 sample data to the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Header row (row 1): columns F..V get new labels; columns A..E stay.
# ---------------------------------------------------------------------
$ws.Range("F1").Value  = "entidad"
$ws.Range("G1").Value  = "monto_alcanza"
$ws.Range("H1").Value  = "plazo"
$ws.Range("I1").Value  = "estado_civil"
$ws.Range("J1").Value  = "tipo_vivienda"
$ws.Range("K1").Value  = "tiempo_pensionado"
$ws.Range("L1").Value  = "contrasena_sipre"
$ws.Range("M1").Value  = "ref1_nombre"
$ws.Range("N1").Value  = "ref1_telefono"
$ws.Range("O1").Value  = "ref1_parentesco"
$ws.Range("P1").Value  = "ref2_nombre"
$ws.Range("Q1").Value  = "ref2_telefono"
$ws.Range("R1").Value  = "ref2_parentesco"
$ws.Range("S1").Value  = "asesor"
$ws.Range("T1").Value  = "asesor_venta"
$ws.Range("U1").Value  = "fuente"
$ws.Range("V1").Value  = "fuente_base_nombre"

# ---------------------------------------------------------------------
# 2) New trailing header columns (W..AC), copying the bold/bordered
#    header style from an existing header cell (A1) so no stray style
#    is introduced.
# ---------------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null
$ws.Range("W1:AC1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("W1").Value  = "fecha_ingreso"
$ws.Range("X1").Value  = "fecha_dispersion"
$ws.Range("Y1").Value  = "fecha_proximo"
$ws.Range("Z1").Value  = "estatus"
$ws.Range("AA1").Value = "observaciones"
$ws.Range("AB1").Value = "score"
$ws.Range("AC1").Value = "sucursal"

# ---------------------------------------------------------------------
# 3) Data row (row 2): fix id value, move values that survive to their
#    new home column, and clear the cells that no longer apply.
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "C1000"
# B2, C2, D2 (nombre, telefono, correo) are unchanged.

$ws.Range("S2").Value  = "Asesor Test"
$ws.Range("U2").Value  = "Test"
$ws.Range("Z2").Value  = "PENDIENTE CLIENTE"
$ws.Range("AA2").Value = "Cliente creado por test de integración"
$ws.Range("AC2").Value = "TOXQUI"

# W2 ("2025-11-27") and AB2 ("700") look like a date / a number, so Excel's
# type-inference would silently store them as a date serial / numeric value.
# The source data keeps them as plain text, so force text storage via the
# "@" number format, then re-normalise the cell format (copy from a plain
# text cell) so no stray number-format style lingers on the cell.
$ws.Range("W2").NumberFormat = "@"
$ws.Range("W2").Value = "2025-11-27"

$ws.Range("AB2").NumberFormat = "@"
$ws.Range("AB2").Value = "700"

$ws.Range("B2").Copy() | Out-Null
$ws.Range("W2").PasteSpecial(-4122) | Out-Null
$ws.Range("AB2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Clear old cells that are no longer part of the populated layout.
$ws.Range("F2").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("Q2").ClearContents()
$ws.Range("R2").ClearContents()
$ws.Range("T2").ClearContents()
